$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A122").Value = "How many curves can be plotted in GEO?"
$ws.Range("B122").Value = "llama3.2:latest"
$ws.Range("C122").Value = "Unfortunately, the provided information does not specify a maximum number of curves that can be plotted in GEO. However, it does mention that a curve can wrap around a maximum of 50 times."

$ws.Range("A123").Value = "What the maximum number of headers I can display in my log?"
$ws.Range("B123").Value = "llama3.2:latest"
$ws.Range("C123").Value = "The maximum number of headers you can display in your log is up to 50."

$ws.Range("A124").Value = "What the maximum number of headers I can display in my log?"
$ws.Range("B124").Value = "llama3.2:latest"
$ws.Range("C124").Value = "The maximum number of headers you can display in your log is up to 50."
